$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop workbook protection (workbookProtection element removed in target).
try { $wb.Unprotect() } catch { }

# New data grid (rows 2-5, columns A-M). Every value in the target XML is
# stored as text (t="inlineStr"/shared-string), even the numeric-looking
# ones ("23", "12.3", "58", ...), so every cell is written as text below.
$data = @(
    @("a1","noxo","lipo+gino","23","12.3","0","58","20","5555","8","9","9","246.0"),
    @("a1","noxo","lipo+gino","23","12.9","0","58","20","5555","8","9","9","258.0"),
    @("a1","noxo","lipo+gino","23","12.9","0","58","13","5555","8","9","9","167.70000000000002"),
    @("v","fr","r","4","12.9","0","0","13","54","5","5","5","167.700")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowIndex = $i + 2
    $row = $data[$i]
    for ($j = 0; $j -lt $row.Length; $j++) {
        $cell = $ws.Cells.Item($rowIndex, $j + 1)
        # Force text storage for numeric-looking tokens (e.g. "23", "58",
        # "246.0") by formatting the cell as Text before assigning, then
        # restore the "Normal" style so no stray style index is left
        # referenced on the cell (matches the target, which carries no
        # custom cell styles).
        $cell.NumberFormat = "@"
        $cell.Value = $row[$j]
        $cell.Style = "Normal"
    }
}

# Match the recorded selection in the target sheet view.
$ws.Range("I8").Select() | Out-Null
